$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q7").Value = 511484
$ws.Range("R7").Value = 6366215
$ws.Range("Z7").ClearContents()
$ws.Range("AB7").ClearContents()
